# Applies the commit's row permutation to rows 2-15 (columns A:AY) of the
# active worksheet. Each full row (all cells) moves to a new row position
# according to the mapping below; row 6 stays in place. This reproduces the
# exact set of cell-level changes described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$minRow = 2
$maxRow = 15
$maxCol = 51   # column AY

# after[r] = before[mapping[r]]
$mapping = @{ 2=15; 3=9; 4=5; 5=2; 6=6; 7=10; 8=13; 9=3; 10=12; 11=4; 12=8; 13=11; 14=7; 15=14 }

# 1) Snapshot every cell value for rows 2..15 BEFORE making any changes.
#    Build each row array with placeholder values first (so index
#    assignment keeps it a proper fixed-length array), then fill it in -
#    this avoids the "+= $null gets swallowed" pitfall.
$data = @{}
for ($r = $minRow; $r -le $maxRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $maxCol; $c++) { $rowVals += 0 }
    for ($c = 1; $c -le $maxCol; $c++) {
        $rowVals[$c-1] = $ws.Cells.Item($r, $c).Value()
    }
    $data[$r] = $rowVals
}

# 2) Write back each destination row using the source row's snapshot.
#    Columns Y (Startdatum) and AA (Slutdatum) hold plain text values that
#    look like dates (e.g. "2023-08-15"); Excel would otherwise silently
#    reinterpret them as date serials on assignment, so force those two
#    columns to keep a text number format before writing them back.
for ($r = $minRow; $r -le $maxRow; $r++) {
    $src = $mapping[$r]
    $srcVals = $data[$src]
    $ws.Cells.Item($r, 25).NumberFormat = "@"
    $ws.Cells.Item($r, 27).NumberFormat = "@"
    for ($c = 1; $c -le $maxCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c-1]
    }
}
